$d = $word.ActiveDocument

# --- 1. Insert a new centered, bold title paragraph at the very top ---
$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.InsertParagraphBefore()

$titlePara = $d.Paragraphs.Item(1)
$titlePara.Alignment = 1   # wdAlignParagraphCenter

$titleText = " Présentation du cours "
$titleRange = $titlePara.Range
$titleRange.Text = $titleText

# Apply bold only to the text run (exclude the paragraph mark) so no
# rPr ends up stamped onto the pPr itself.
$boldRange = $d.Range($titleRange.Start, $titleRange.Start + $titleText.Length)
$boldRange.Font.Bold = 1

# --- 2. Center the "fin de la liste qui n'est pas vide" paragraph ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "fin de la liste qui n'est pas vide*") {
        $p.Alignment = 1   # wdAlignParagraphCenter
    }
}
